$d = $word.ActiveDocument

$replacements = @(
    @{old='549×3=1647'; new='995×2=1990'},
    @{old='849×9=7641'; new='243×5=1215'},
    @{old='285×9=2565'; new='646×8=5168'},
    @{old='456×8=3648'; new='768×8=6144'},
    @{old='561×4=2244'; new='680×9=6120'},
    @{old='670×5=3350'; new='556×7=3892'},
    @{old='800×8=6400'; new='433×5=2165'},
    @{old='283×6=1698'; new='612×2=1224'},
    @{old='945×3=2835'; new='241×8=1928'},
    @{old='880×6=5280'; new='957×7=6699'},
    @{old='551×9=4959'; new='730×6=4380'},
    @{old='221×6=1326'; new='576×6=3456'},
    @{old='715×8=5720'; new='248×7=1736'},
    @{old='152×8=1216'; new='907×8=7256'},
    @{old='655×7=4585'; new='228×4=912'},
    @{old='106×3=318'; new='747×6=4482'},
    @{old='509×7=3563'; new='913×7=6391'},
    @{old='551×2=1102'; new='320×6=1920'},
    @{old='527×9=4743'; new='866×3=2598'},
    @{old='361×6=2166'; new='935×7=6545'},
    @{old='584×2=1168'; new='212×3=636'},
    @{old='396×7=2772'; new='139×6=834'},
    @{old='232×2=464'; new='641×5=3205'},
    @{old='884×2=1768'; new='870×4=3480'},
    @{old='526×6=3156'; new='319×4=1276'}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
